$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.685.62"
$ws.Range("E2").Value = "  +4.93%  "
$ws.Range("D3").Value = "2.458.95"
$ws.Range("E3").Value = "  +3.04%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "159.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "496.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.71%  "
$ws.Range("E7").Value = "  +23.54%  "
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").Value = "2.484.13"
$ws.Range("E9").Value = "  +3.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +15.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.338"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.36%  "
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "2.881.69"
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").Value = "58.554.88"
$ws.Range("E15").Value = "  +3.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").Value = "2.482.13"
$ws.Range("E18").Value = "  +3.91%  "
$ws.Range("E19").Value = "  +5.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "329.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.18%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.411"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.86%  "
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "2.567.74"
$ws.Range("E29").Value = "  +2.14%  "
$ws.Range("D30").Value = "0.0₃0809"
$ws.Range("E30").Value = "  +5.00%  "
$ws.Range("E31").Value = "  -0.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "152.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.32%  "
$ws.Range("E34").Value = "  +3.85%  "
$ws.Range("E35").Value = "  +8.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.99%  "
$ws.Range("E37").Value = "  +5.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.848"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.41%  "
$ws.Range("E40").Value = "  +6.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "34.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.08%  "
$ws.Range("E43").Value = "  +7.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.609"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.991"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("E47").Value = "  +5.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.708"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.48%  "
